$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15: add date and activity text
$ws.Range("A15").Value = (Get-Date -Year 2015 -Month 8 -Day 12 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B15").Value = "Implemented clpp radix sort (not working), Implemented CPU Scan"

# Update the active selection to B27
$ws.Activate()
$ws.Range("B27").Select()
